$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new-string-introducing cells first, in the exact order they were
# originally typed, so shared strings land in the same order as the source.
$ws.Range("A23").Value = "Выдать действующие этажи"
$ws.Range("A24").Value = "Выдать все этажи"
$ws.Range("D23").Value = "Номера и названия этажей"
$ws.Range("D24").Value = "Номера и названия этажей"
$ws.Range("A27").Value = "Загрузить схему"
$ws.Range("A21").Value = "Выдать id следующей точки"
$ws.Range("A22").Value = "Выдать id следующего Этажа"
$ws.Range("A25").Value = "Создать менеджера"
$ws.Range("A26").Value = "Загрузить список всех менеджеров"
$ws.Range("D27").Value = "Схема"
$ws.Range("D26").Value = "логины всех менеджеров"
$ws.Range("D21").Value = "Свободный номер точки"
$ws.Range("D22").Value = "Свободный номер этажа"
$ws.Range("C27").Value = "id"

# Remaining cells (reuse already-existing shared strings / plain numbers)
$ws.Range("B21").Value = 17
$ws.Range("C21").Value = "-"

$ws.Range("B22").Value = 18
$ws.Range("C22").Value = "-"

$ws.Range("B23").Value = 19
$ws.Range("C23").Value = "-"

$ws.Range("B24").Value = 20
$ws.Range("C24").Value = "-"

$ws.Range("B25").Value = 21
$ws.Range("C25").Value = "Логин + Пароль"
$ws.Range("D25").Value = "true/false"

$ws.Range("B26").Value = 22
$ws.Range("C26").Value = "-"

$ws.Range("B27").Value = 23

# Restore center-alignment style on columns B and C for the new rows, matching
# the formatting already used throughout the rest of the table.
$ws.Range("B21:C27").HorizontalAlignment = -4108  # xlCenter

# Update the saved view/selection to match the edited state.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("C28").Select()
